$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Expand the student name list in the first paragraph.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Student Name: Yinan Du", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Student Name: Yinan Du, Zilong Zheng, Brian Luu", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Re-join the sentence that the old "_GoBack" bookmark used to split into
#    "...Price To Sales Ratio, Discoun" / "ted Cash Flow..." two runs. Doing a
#    find/replace over the full (unchanged) sentence text removes the stray
#    bookmark that sat inside it and merges the text back into a single run.
# ---------------------------------------------------------------------------
$fullSentence = "Price To Sales Ratio, Discounted Cash Flow. Information will be based on date 12/03/2021."
$d.Content.Find.Execute($fullSentence, $true, $false, $false, $false, $false,
                         $true, 1, $false, $fullSentence, 2) | Out-Null

# Safety net: if for any reason the bookmark survived the replace above, remove it
# explicitly so the only "_GoBack" bookmark left is the one (re)created in step 3.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 3. Re-create the "_GoBack" bookmark, collapsed, at the end of the first
#    paragraph (immediately after the updated name list, before the
#    paragraph mark). A bookmark can't be added directly at that boundary
#    position in one step, so a temporary placeholder character is inserted
#    there first, the bookmark is anchored just before it, and the
#    placeholder is removed again.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range
$endPos = $p1.End - 1
$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("X")
$bookmarkPoint = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bookmarkPoint)
$d.Range($endPos, $endPos + 1).Delete()
